$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45972
$ws.Range("B2").Value = 91.28
$ws.Range("C2").Value = 88.65000000000001
$ws.Range("D2").Value = 89.22
$ws.Range("E2").Value = 88.95
$ws.Range("F2").Value = 88.70999999999999
$ws.Range("G2").Value = 88.44
$ws.Range("H2").Value = 89.66
$ws.Range("I2").Value = 94.87
$ws.Range("J2").Value = 93.43000000000001
$ws.Range("K2").Value = 50.99
$ws.Range("L2").Value = 13.07
$ws.Range("M2").Value = 10.7
$ws.Range("N2").Value = 12.96
$ws.Range("O2").Value = 8.67
$ws.Range("P2").Value = 5.64
$ws.Range("Q2").Value = 12.82
$ws.Range("R2").Value = 42.01
$ws.Range("S2").Value = 79.66
$ws.Range("T2").Value = 95.78
$ws.Range("U2").Value = 102.48
$ws.Range("V2").Value = 102.61
$ws.Range("W2").Value = 96.17
$ws.Range("X2").Value = 87.89
$ws.Range("Y2").Value = 70.45999999999999
$ws.Range("Z2").Value = 66.88
$ws.Range("AA2").Value = "4h-8h"
$ws.Range("AB2").Value = 90.42
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 99.39
$ws.Range("AE2").Value = "18h-20h"
$ws.Range("AF2").Value = 99.13
$ws.Range("AG2").Value = "9h-16h"
